$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the merged requisito text in row 23 (B23/C23)
$newText = "LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)`n"
$ws.Range("B23").Value = $newText
$ws.Range("C23").Value = $newText

# Delete row 24 entirely (was LOM3221 requisito row)
$ws.Rows("24").Delete()
